$d = $word.ActiveDocument

$replacements = @(
    @("377×4=1508", "793×2=1586"),
    @("754×8=6032", "190×7=1330"),
    @("830×2=1660", "692×7=4844"),
    @("429×3=1287", "450×3=1350"),
    @("905×7=6335", "905×2=1810"),
    @("151×9=1359", "595×2=1190"),
    @("223×9=2007", "473×7=3311"),
    @("222×4=888",  "442×5=2210"),
    @("951×5=4755", "338×3=1014"),
    @("188×7=1316", "471×2=942"),
    @("726×2=1452", "767×7=5369"),
    @("381×5=1905", "672×2=1344"),
    @("228×9=2052", "186×9=1674"),
    @("735×9=6615", "663×3=1989"),
    @("309×3=927",  "266×8=2128"),
    @("304×2=608",  "591×8=4728"),
    @("423×8=3384", "862×3=2586"),
    @("838×7=5866", "378×6=2268"),
    @("832×5=4160", "921×6=5526"),
    @("506×4=2024", "290×6=1740"),
    @("434×7=3038", "512×7=3584"),
    @("488×2=976",  "138×7=966"),
    @("356×7=2492", "850×6=5100"),
    @("552×5=2760", "313×2=626"),
    @("894×4=3576", "986×2=1972")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $find = $d.Content.Find
    $find.ClearFormatting()
    $find.Replacement.ClearFormatting()
    $find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
